# The edit reshuffles the data rows (2..26) of the single sheet: each target
# row ends up with the full contents (all columns A:T) that some other row
# used to hold, i.e. the whole data block is permuted row-for-row (dates,
# quality grade, volumes, prices, units, origin, etc. all travel together).
#
# Mapping: target row -> source row (the row whose original contents now
# live in the target row).
$map = @{
    2  = 18
    3  = 17
    4  = 9
    5  = 10
    6  = 13
    7  = 21
    8  = 7
    9  = 20
    10 = 6
    11 = 26
    12 = 19
    13 = 24
    14 = 25
    15 = 8
    16 = 3
    17 = 11
    18 = 12
    19 = 22
    20 = 2
    21 = 15
    22 = 14
    23 = 5
    24 = 4
    25 = 16
    26 = 23
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every source row's full A:T contents BEFORE writing anything back,
# since rows are both sources and destinations (a permutation).
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":T" + $r).Value2
}

foreach ($target in $map.Keys) {
    $source = $map[$target]
    $ws.Range("A" + $target + ":T" + $target).Value = $snapshot[$source]
}
